# Generate Report for handoff
# Adds two new localized-file rows (2ca66774-... and 6f4488dc-...) to the
# Overview / zh-cn / de-de sheets, ahead of the existing ".localization-config"
# row, and flips the status of the two pre-existing in-flight files
# (0356ee0b-... and 1fe862d5-...) from "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2        # xlUnderlineStyleSingle
$hyperlinkColor     = 15570276 # RGB(0x64,0x95,0xED) == "FF6495ED"
$dateFormat         = "yyyy-mm-dd HH:mm:ss"

function Style-AsName([object]$rng) {
    $rng.Font.Underline = $hyperlinkUnderline
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Existing rows 2 & 3 (0356ee0b.md / 1fe862d5.md) move from
# "Ready for handoff" to "In Translation"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"

# Row 4 becomes the new 2ca66774 file (was ".localization-config")
$wsOverview.Range("A4").Value = "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md", "", "", "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md")
Style-AsName($wsOverview.Range("A4"))

# Row 5 (new) is the 6f4488dc file
$wsOverview.Range("A5").Value = "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md", "", "", "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md")
Style-AsName($wsOverview.Range("A5"))

# Row 6 (new) is the relocated ".localization-config" row
$wsOverview.Range("A6").Value = ".localization-config"
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/.localization-config", "", "", ".localization-config")
Style-AsName($wsOverview.Range("A6"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "In Translation"
$wsZh.Range("B3").Value = "In Translation"

# Row 4 becomes the new 2ca66774 file data (was ".localization-config")
$wsZh.Range("A4").Value = "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-01-19 04:03:56"
$wsZh.Range("D4").NumberFormat = $dateFormat
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md", "", "", "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md")
Style-AsName($wsZh.Range("A4"))
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c55623bb9b3278aea10b519147b5e790d27d6695/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.zh-cn.xlf", "", "", "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.zh-cn.xlf")
Style-AsName($wsZh.Range("C4"))

# Row 5 (new) is the 6f4488dc file data
$wsZh.Range("A5").Value = "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.zh-cn.xlf"
$wsZh.Range("D5").Value = "2016-01-19 04:03:56"
$wsZh.Range("D5").NumberFormat = $dateFormat
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md", "", "", "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md")
Style-AsName($wsZh.Range("A5"))
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c55623bb9b3278aea10b519147b5e790d27d6695/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.zh-cn.xlf", "", "", "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.zh-cn.xlf")
Style-AsName($wsZh.Range("C5"))

# Row 6 (new) is the relocated ".localization-config" row (no C cell)
$wsZh.Range("A6").Value = ".localization-config"
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("D6").NumberFormat = $dateFormat
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/.localization-config", "", "", ".localization-config")
Style-AsName($wsZh.Range("A6"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "In Translation"
$wsDe.Range("B3").Value = "In Translation"

# Row 4 becomes the new 2ca66774 file data (was ".localization-config")
$wsDe.Range("A4").Value = "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.de-de.xlf"
$wsDe.Range("D4").Value = "2016-01-19 04:04:06"
$wsDe.Range("D4").NumberFormat = $dateFormat
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md", "", "", "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.md")
Style-AsName($wsDe.Range("A4"))
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b17aadf3e4f6b83ff95f7bbad9761ecaab0310a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.de-de.xlf", "", "", "2ca66774-1a61-48dd-98ea-6413ab8b7a7e.d42b24457daf4892b673cd6678339a0afc014ea0.de-de.xlf")
Style-AsName($wsDe.Range("C4"))

# Row 5 (new) is the 6f4488dc file data
$wsDe.Range("A5").Value = "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.de-de.xlf"
$wsDe.Range("D5").Value = "2016-01-19 04:04:06"
$wsDe.Range("D5").NumberFormat = $dateFormat
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/e2e/6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md", "", "", "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.md")
Style-AsName($wsDe.Range("A5"))
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b17aadf3e4f6b83ff95f7bbad9761ecaab0310a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.de-de.xlf", "", "", "6f4488dc-e4f8-48cd-8f32-ad80637efa3d.850a94636b5e844744ad80cd3fb7e43966ba3d8d.de-de.xlf")
Style-AsName($wsDe.Range("C5"))

# Row 6 (new) is the relocated ".localization-config" row (no C cell)
$wsDe.Range("A6").Value = ".localization-config"
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("D6").NumberFormat = $dateFormat
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/94ade2b7c1cdbd4d5c7d502d02002a8dc2b4a5cf/.localization-config", "", "", ".localization-config")
Style-AsName($wsDe.Range("A6"))
